# Fix casing/typos in the class-name labels (column A) so the site keys
# consistently use camelCase naming, matching the dropdown / non-dropdown
# class list used across all sites.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value2 = "mdaTextHomePage"
$ws.Range("A4").Value2 = "mdaTitle"
$ws.Range("A8").Value2 = "pageTitleNewTab"

# Move the active selection to A8, matching the sheet view state saved
# with the workbook after the edit.
$ws.Range("A8").Select()
